# Remove the specific rows that were deleted upstream, causing all
# subsequent rows to shift up. Row numbers refer to the *original*
# worksheet (1-based, including the header row), so we delete them in
# descending order to keep earlier row numbers valid.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PEBCOM")

$rowsToDelete = @(49, 47, 46, 39, 33, 23, 10, 6)

foreach ($r in $rowsToDelete) {
    $ws.Rows($r).Delete()
}
